# Auto-generated script applying numeric updates described in the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 151.71428
$ws.Range("I19").Value = 151.71428
$ws.Range("K19").Value = 151.71428
$ws.Range("M19").Value = 23.28572
$ws.Range("H86").Value = 2358.2
$ws.Range("I86").Value = 1372.75
$ws.Range("K86").Value = 1372.75
$ws.Range("M86").Value = -249.75
$ws.Range("H89").Value = 2358.2
$ws.Range("I89").Value = 1372.75
$ws.Range("K89").Value = 6863.75
$ws.Range("M89").Value = -1247.75
$ws.Range("H92").Value = 584
$ws.Range("I92").Value = 625.75
$ws.Range("K92").Value = 625.75
$ws.Range("M92").Value = 622.25
$ws.Range("H98").Value = 3821.7334
$ws.Range("I98").Value = 3650.111
$ws.Range("K98").Value = 3650.111
$ws.Range("M98").Value = -2152.111
$ws.Range("H116").Value = 6308.5835
$ws.Range("I116").Value = 5766.2856
$ws.Range("J116").Value = 7067.8
$ws.Range("K116").Value = 5766.2856
$ws.Range("L116").Value = 7067.8
$ws.Range("M116").Value = -2324.2856
$ws.Range("N116").Value = -13951.8
$ws.Range("H122").Value = 3821.7334
$ws.Range("I122").Value = 3650.111
$ws.Range("K122").Value = 10950.333
$ws.Range("M122").Value = -8500.332999999999
$ws.Range("H129").Value = 3502.9412
$ws.Range("I129").Value = 1017.2
$ws.Range("K129").Value = 3051.6
$ws.Range("M129").Value = 1948.4
$ws.Range("H132").Value = 2904.756
$ws.Range("I132").Value = 2736.3242
$ws.Range("K132").Value = 8208.972600000001
$ws.Range("M132").Value = -5678.972600000001
$ws.Range("H138").Value = 2504.8
$ws.Range("I138").Value = 995.3333
$ws.Range("K138").Value = 2985.9999
$ws.Range("M138").Value = 2154.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5541.6665
$ws.Range("I32").Value = 4531.25
$ws.Range("K32").Value = 4531.25
$ws.Range("M32").Value = -4244.25
$ws.Range("H45").Value = 3554.7
$ws.Range("I45").Value = 3749.4
$ws.Range("K45").Value = 3749.4
$ws.Range("M45").Value = -3372.4
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("H74").Value = 1373.4615
$ws.Range("I74").Value = 1150.4546
$ws.Range("K74").Value = 1150.4546
$ws.Range("M74").Value = -276.4546
$ws.Range("H77").Value = 1373.4615
$ws.Range("I77").Value = 1150.4546
$ws.Range("K77").Value = 5752.273
$ws.Range("M77").Value = -1384.273
$ws.Range("H102").Value = 2500
$ws.Range("I102").Value = 2500
$ws.Range("K102").Value = 2500
$ws.Range("M102").Value = -878
$ws.Range("H132").Value = 5710.364
$ws.Range("I132").Value = 5180
$ws.Range("K132").Value = 15540
$ws.Range("M132").Value = -13010

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1475.9
$ws.Range("I94").Value = 1466.6786
$ws.Range("J94").Value = 1605
$ws.Range("K94").Value = 1466.6786
$ws.Range("L94").Value = 1605
$ws.Range("M94").Value = -1015.6786
$ws.Range("N94").Value = -2507
$ws.Range("H107").Value = 1948.8572
$ws.Range("I107").Value = 910.5
$ws.Range("K107").Value = 910.5
$ws.Range("M107").Value = 1009.5
$ws.Range("H127").Value = 31998
$ws.Range("J127").Value = 31998
$ws.Range("L127").Value = 31998
$ws.Range("N127").Value = -41918
$ws.Range("H134").Value = 7637.375
$ws.Range("I134").Value = 7637.375
$ws.Range("K134").Value = 22912.125
$ws.Range("M134").Value = -20377.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 295
$ws.Range("I4").Value = 295
$ws.Range("K4").Value = 295
$ws.Range("M4").Value = -183
$ws.Range("H31").Value = 3388.7222
$ws.Range("I31").Value = 2226.8333
$ws.Range("J31").Value = 5712.5
$ws.Range("K31").Value = 2226.8333
$ws.Range("L31").Value = 5712.5
$ws.Range("M31").Value = -1931.8333
$ws.Range("N31").Value = -6302.5
$ws.Range("H34").Value = 3388.7222
$ws.Range("I34").Value = 2226.8333
$ws.Range("J34").Value = 5712.5
$ws.Range("K34").Value = 2226.8333
$ws.Range("L34").Value = 5712.5
$ws.Range("M34").Value = -2024.8333
$ws.Range("N34").Value = -6116.5
$ws.Range("H122").Value = 1529.8334
$ws.Range("I122").Value = 1495.9
$ws.Range("J122").Value = 1699.5
$ws.Range("K122").Value = 4487.700000000001
$ws.Range("L122").Value = 5098.5
$ws.Range("M122").Value = -2037.700000000001
$ws.Range("N122").Value = -9998.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1905.2222
$ws.Range("I107").Value = 1774.25
$ws.Range("J107").Value = 2010
$ws.Range("K107").Value = 5322.75
$ws.Range("L107").Value = 6030
$ws.Range("M107").Value = -3402.75
$ws.Range("N107").Value = -9870
$ws.Range("H122").Value = 1589.4615
$ws.Range("I122").Value = 1529.3334
$ws.Range("J122").Value = 1641
$ws.Range("K122").Value = 13764.0006
$ws.Range("L122").Value = 14769
$ws.Range("M122").Value = -11314.0006
$ws.Range("N122").Value = -19669
$ws.Range("H129").Value = 1973.8572
$ws.Range("I129").Value = 971.6
$ws.Range("J129").Value = 2530.6667
$ws.Range("K129").Value = 2914.8
$ws.Range("L129").Value = 7592.000100000001
$ws.Range("M129").Value = 2085.2
$ws.Range("N129").Value = -17592.0001
$ws.Range("H131").Value = 2393.125
$ws.Range("I131").Value = 1579
$ws.Range("J131").Value = 2763.182
$ws.Range("K131").Value = 4737
$ws.Range("L131").Value = 8289.545999999998
$ws.Range("M131").Value = 303
$ws.Range("N131").Value = -18369.546

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5000
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 5000
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 5000
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -6996
$ws.Range("H83").Value = 5000
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 5000
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 25000
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -34984
$ws.Range("H107").Value = 255.6
$ws.Range("I107").Value = 255.6
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 255.6
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1664.4
$ws.Range("N107").ClearContents()
$ws.Range("H132").Value = 793
$ws.Range("I132").Value = 793
$ws.Range("K132").Value = 2379
$ws.Range("M132").Value = 151

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 12440.357
$ws.Range("I16").Value = 1925.2727
$ws.Range("J16").Value = 50995.668
$ws.Range("K16").Value = 1925.2727
$ws.Range("L16").Value = 50995.668
$ws.Range("M16").Value = -1755.2727
$ws.Range("N16").Value = -51335.668
$ws.Range("H55").Value = 1464
$ws.Range("I55").Value = 1840
$ws.Range("J55").Value = 900
$ws.Range("K55").Value = 1840
$ws.Range("L55").Value = 900
$ws.Range("M55").Value = -1667
$ws.Range("N55").Value = -1246
$ws.Range("H122").Value = 5301.6
$ws.Range("I122").Value = 3502.6667
$ws.Range("J122").Value = 8000
$ws.Range("K122").Value = 10508.0001
$ws.Range("L122").Value = 24000
$ws.Range("M122").Value = -8058.000100000001
$ws.Range("N122").Value = -28900
$ws.Range("H136").Value = 3664.1052
$ws.Range("I136").Value = 3301.4
$ws.Range("K136").Value = 9904.200000000001
$ws.Range("M136").Value = -7354.200000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 2500
$ws.Range("J18").Value = 2500
$ws.Range("L18").Value = 2500
$ws.Range("N18").Value = -2846
$ws.Range("H51").Value = 28499.5
$ws.Range("I51").Value = 28499.5
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 28499.5
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -27989.5
$ws.Range("N51").ClearContents()
$ws.Range("H96").Value = 1334.3334
$ws.Range("I96").Value = 999
$ws.Range("J96").Value = 1502
$ws.Range("K96").Value = 999
$ws.Range("L96").Value = 1502
$ws.Range("M96").Value = 374
$ws.Range("N96").Value = -4248
$ws.Range("H100").Value = 809.9167
$ws.Range("I100").Value = 777.75
$ws.Range("J100").Value = 874.25
$ws.Range("K100").Value = 1555.5
$ws.Range("L100").Value = 1748.5
$ws.Range("M100").Value = -1014.5
$ws.Range("N100").Value = -2830.5
$ws.Range("H111").Value = 37660.75
$ws.Range("J111").Value = 37660.75
$ws.Range("L111").Value = 37660.75
$ws.Range("N111").Value = -45840.75
$ws.Range("H132").Value = 2499
$ws.Range("I132").Value = 2499
$ws.Range("K132").Value = 7497
$ws.Range("M132").Value = -4967
